$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-27 Thursday", "2025-03-28 Friday"),
    @("72÷5=", "50÷7="),
    @("12÷7=", "57÷8="),
    @("17÷9=", "74÷5="),
    @("63÷2=", "66÷3="),
    @("57÷2=", "73÷2="),
    @("62÷7=", "80÷9="),
    @("76÷6=", "67÷8="),
    @("95÷4=", "35÷6="),
    @("33÷3=", "75÷9="),
    @("67÷4=", "75÷7="),
    @("96÷6=", "89÷4="),
    @("38÷2=", "47÷4="),
    @("78÷2=", "57÷7="),
    @("45÷7=", "74÷8="),
    @("97÷3=", "34÷7="),
    @("28÷5=", "69÷7="),
    @("72÷8=", "64÷2="),
    @("13÷5=", "95÷3="),
    @("44÷6=", "48÷3="),
    @("65÷3=", "17÷2="),
    @("60÷2=", "50÷4="),
    @("60÷7=", "30÷9="),
    @("52÷9=", "68÷4="),
    @("47÷5=", "61÷4="),
    @("42÷6=", "57÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
